$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Suspended Sediment Concentration" was dropped from both the 5-year and
# 10-year trend-period blocks. Deleting row 11 shifts every later row up by
# one, so the second "Suspended Sediment Concentration" row -- originally
# row 24 -- is now at row 23.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(23).Delete()

# Refresh the recalculated trend-analysis figures for the remaining rows.

# Row 2
$ws.Range("E2").Value = "WARNING: Sen slope based on tied non-censored values"
$ws.Range("F2").Value = 0.40655411136905
$ws.Range("G2").Value = 0.024390243902439
$ws.Range("H2").Value = 0.634146341463415
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = -0.242691029900332
$ws.Range("M2").Value = 0.164386613265134
$ws.Range("N2").Value = 0
$ws.Range("P2").Value = "As likely as not improving"

# Row 3
$ws.Range("F3").Value = 0.102241480749488
$ws.Range("H3").Value = 0.945454545454545
$ws.Range("J3").Value = 11.31
$ws.Range("K3").Value = -0.0880918003367599
$ws.Range("L3").Value = -0.143931743602562
$ws.Range("M3").Value = 0.006716199912325
$ws.Range("N3").Value = -0.778884176275508

# Row 4
$ws.Range("F4").Value = 0.460427237813098
$ws.Range("H4").Value = 0.357142857142857
$ws.Range("J4").Value = 0.008
$ws.Range("L4").Value = -0.0006819258417971
$ws.Range("M4").Value = 0.0007435159590956

# Row 5
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "ok"
$ws.Range("F5").Value = 0.152036016958434
$ws.Range("G5").Value = 0.0363636363636364
$ws.Range("H5").Value = 0.672727272727273
$ws.Range("J5").Value = 49
$ws.Range("K5").Value = 5.98770491803279
$ws.Range("L5").Value = -4.19957280099892
$ws.Range("M5").Value = 13.1895796128623
$ws.Range("N5").Value = 12.219805955169
$ws.Range("P5").Value = "Unlikely improving"

# Row 6
$ws.Range("F6").Value = 0.168575373327403
$ws.Range("G6").Value = 0.867924528301887
$ws.Range("H6").Value = 0.188679245283019
$ws.Range("I6").Value = 4

# Row 7
$ws.Range("F7").Value = 0.918578376120242
$ws.Range("P7").Value = "Very likely improving"

# Row 8
$ws.Range("F8").Value = 0.0002815952365126
$ws.Range("G8").Value = 0.232142857142857
$ws.Range("J8").Value = 0.044
$ws.Range("K8").Value = 0.0099192972632463
$ws.Range("L8").Value = 0.0033495438815968
$ws.Range("M8").Value = 0.0149038633087759
$ws.Range("N8").Value = 22.5438574164689
$ws.Range("P8").Value = "Exceptionally unlikely improving"

# Row 9
$ws.Range("D9").Value = $true
$ws.Range("F9").Value = 0.019607004146115
$ws.Range("J9").Value = 7.84
$ws.Range("K9").Value = -0.0569974464463047
$ws.Range("L9").Value = -0.116722655839948
$ws.Range("M9").Value = -0.025837030269252
$ws.Range("N9").Value = -0.72700824548858
$ws.Range("P9").Value = "Extremely unlikely increasing"

# Row 10
$ws.Range("F10").Value = 0.0060978008349758
$ws.Range("G10").Value = 0.0714285714285714
$ws.Range("H10").Value = 0.839285714285714
$ws.Range("J10").Value = 0.055
$ws.Range("K10").Value = 0.0133046372626703
$ws.Range("L10").Value = 0.0028365096517086
$ws.Range("M10").Value = 0.017527517437631
$ws.Range("N10").Value = 24.1902495684914
$ws.Range("P10").Value = "Exceptionally unlikely improving"

# Row 11
$ws.Range("B11").Value = "Total Nitrogen"
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = "ok"
$ws.Range("F11").Value = 0.0214800730142125
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0.16
$ws.Range("K11").Value = 0.0183007260596546
$ws.Range("L11").Value = 0.00405268069938
$ws.Range("M11").Value = 0.0277024917936252
$ws.Range("N11").Value = 11.4379537872841
$ws.Range("P11").Value = "Extremely unlikely improving"
$ws.Range("W11").Value = "g/m3"

# Row 12
$ws.Range("B12").Value = "Total Phosphorus"
$ws.Range("D12").Value = $false
$ws.Range("E12").Value = "ok"
$ws.Range("F12").Value = 0.112475455343972
$ws.Range("H12").Value = 0.446428571428571
$ws.Range("J12").Value = 0.015
$ws.Range("K12").Value = 0.0006523340700555
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0.0018755204596758
$ws.Range("N12").Value = 4.34889380037059
$ws.Range("P12").Value = "Unlikely improving"

# Row 13
$ws.Range("B13").Value = "Turbidity"
$ws.Range("D13").Value = $true
$ws.Range("F13").Value = 0.5
$ws.Range("H13").Value = 0.857142857142857
$ws.Range("J13").Value = 1.04
$ws.Range("K13").Value = -0.0057560626873126
$ws.Range("L13").Value = -0.116391218729894
$ws.Range("M13").Value = 0.136477347730091
$ws.Range("N13").Value = -0.553467566087753
$ws.Range("P13").Value = "As likely as not improving"
$ws.Range("W13").Value = "NTU/FNU"

# Row 14
$ws.Range("B14").Value = "Visual Clarity"
$ws.Range("C14").Value = 10
$ws.Range("E14").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F14").Value = 0.022837762713854
$ws.Range("G14").Value = 0.0842105263157895
$ws.Range("H14").Value = 0.557894736842105
$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 2.75
$ws.Range("K14").Value = -0.06646951774340321
$ws.Range("L14").Value = -0.145812515486057
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -2.41707337248739
$ws.Range("P14").Value = "Extremely unlikely improving"
$ws.Range("W14").Value = "m"

# Row 15
$ws.Range("B15").Value = "Dissolved Oxygen Concentration"
$ws.Range("D15").Value = $true
$ws.Range("E15").Value = "ok"
$ws.Range("F15").Value = 0.945187234608619
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0.876106194690266
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 11.06
$ws.Range("K15").Value = 0.0262885669312751
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0.0599924743367947
$ws.Range("N15").Value = 0.237690478582958
$ws.Range("P15").Value = "Very likely increasing"
$ws.Range("W15").Value = "g/m3"

# Row 16
$ws.Range("B16").Value = "Dissolved Reactive Phosphorus"
$ws.Range("D16").Value = $false
$ws.Range("F16").Value = 0.999852177095218
$ws.Range("H16").Value = 0.258620689655172
$ws.Range("J16").Value = 0.01
$ws.Range("K16").Value = -0.0006156817164735
$ws.Range("L16").Value = -0.000871285660414
$ws.Range("M16").Value = -0.0003519975556925
$ws.Range("N16").Value = -6.15681716473572
$ws.Range("P16").Value = "Virtually certain improving"
$ws.Range("W16").Value = "mg/L"

# Row 17
$ws.Range("B17").Value = "E. coli"
$ws.Range("F17").Value = 0.07339832358813279
$ws.Range("G17").Value = 0.0263157894736842
$ws.Range("H17").Value = 0.491228070175439
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 48
$ws.Range("K17").Value = 2.24692412850308
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 5.48197879705584
$ws.Range("N17").Value = 4.68109193438141
$ws.Range("P17").Value = "Very unlikely improving"
$ws.Range("W17").Value = "E. coli/100 mL"

# Row 18
$ws.Range("B18").Value = "Ammoniacal Nitrogen (NH4)"
$ws.Range("F18").Value = 0.231934062630747
$ws.Range("G18").Value = 0.905660377358491
$ws.Range("H18").Value = 0.216981132075472
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 0.005
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("P18").Value = "Unlikely improving"
$ws.Range("W18").Value = "mg/L"

# Row 19
$ws.Range("B19").Value = "Nitrite Nitrogen (NO2)"
$ws.Range("F19").Value = 0.9997183811561799
$ws.Range("G19").Value = 0.594827586206897
$ws.Range("H19").Value = 0.0862068965517241
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 0.001
$ws.Range("P19").Value = "Virtually certain improving"

# Row 20
$ws.Range("B20").Value = "Nitrate Nitrogen (NO3)"
$ws.Range("D20").Value = $true
$ws.Range("E20").Value = "ok"
$ws.Range("F20").Value = 0.007455847251571
$ws.Range("G20").Value = 0.275862068965517
$ws.Range("H20").Value = 0.629310344827586
$ws.Range("J20").Value = 0.0435
$ws.Range("K20").Value = 0.000456877217169
$ws.Range("M20").Value = 0.0021855689277995
$ws.Range("N20").Value = 1.05029245326226
$ws.Range("P20").Value = "Exceptionally unlikely improving"

# Row 21
$ws.Range("B21").Value = "pH"
$ws.Range("E21").Value = "ok"
$ws.Range("F21").Value = 0.361338159310276
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0.699115044247788
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 7.89
$ws.Range("K21").Value = -0.0048068366326344
$ws.Range("L21").Value = -0.0245867614061331
$ws.Range("M21").Value = 0.0162137243549169
$ws.Range("N21").Value = -0.0609231512374456
$ws.Range("P21").Value = "As likely as not increasing"

# Row 22
$ws.Range("B22").Value = "SIN (Soluble Inorganic nitrogen)"
$ws.Range("D22").Value = $true
$ws.Range("F22").Value = 0.07307266465746989
$ws.Range("G22").Value = 0.103448275862069
$ws.Range("H22").Value = 0.681034482758621
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 0.049
$ws.Range("K22").Value = 0.0012921326666171
$ws.Range("L22").Value = -0.0000119541347115276
$ws.Range("M22").Value = 0.0031250855900789
$ws.Range("N22").Value = 2.63700544207583
$ws.Range("P22").Value = "Very unlikely improving"
$ws.Range("W22").Value = "g/m3"

# Row 23
$ws.Range("B23").Value = "Total Nitrogen"
$ws.Range("E23").Value = "ok"
$ws.Range("F23").Value = 0.231657353217304
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0.431034482758621
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0.15
$ws.Range("K23").Value = 0.0025819287908997
$ws.Range("L23").Value = -0.0016685298353537
$ws.Range("M23").Value = 0.0071477495107632
$ws.Range("N23").Value = 1.7212858605998
$ws.Range("P23").Value = "Unlikely improving"

# Row 24
$ws.Range("B24").Value = "Total Phosphorus"
$ws.Range("E24").Value = "ok"
$ws.Range("F24").Value = 0.121413250282443
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0.267241379310345
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0.014
$ws.Range("K24").Value = 0.0001713586604831
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0.0005180088363264
$ws.Range("N24").Value = 1.22399043202268
$ws.Range("P24").Value = "Unlikely improving"
$ws.Range("W24").Value = "g/m3"

# Row 25
$ws.Range("B25").Value = "Turbidity"
$ws.Range("D25").Value = $false
$ws.Range("F25").Value = 0.288289819896242
$ws.Range("H25").Value = 0.827586206896552
$ws.Range("J25").Value = 1.15
$ws.Range("K25").Value = 0.013091566032888
$ws.Range("L25").Value = -0.0275401235382334
$ws.Range("M25").Value = 0.0546609154524742
$ws.Range("N25").Value = 1.13839704633808
$ws.Range("P25").Value = "Unlikely improving"
$ws.Range("W25").Value = "NTU/FNU"
